# Insert 4 new weekly rows for Chirimoya "Cultivar IV Región" (week of 44466)
# above the existing block, shifting existing rows 53-99 down to 57-103.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("53:56").Insert()

$ws.Cells.Item(53,1).Value = 6
$ws.Cells.Item(53,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(53,3).Value = "Metropolitana"
$ws.Cells.Item(53,4).Value = 44466
$ws.Cells.Item(53,5).Value = 13
$ws.Cells.Item(53,6).Value = "Fruta"
$ws.Cells.Item(53,7).Value = 100107
$ws.Cells.Item(53,8).Value = "Otros"
$ws.Cells.Item(53,9).Value = 100107002
$ws.Cells.Item(53,10).Value = "Chirimoya"
$ws.Cells.Item(53,11).Value = "Cultivar IV Región"
$ws.Cells.Item(53,12).Value = "Especial"
$ws.Cells.Item(53,13).Value = 175
$ws.Cells.Item(53,14).Value = 3100
$ws.Cells.Item(53,15).Value = 3100
$ws.Cells.Item(53,16).Value = 3100
$ws.Cells.Item(53,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(53,18).Value = "Provincia de Limarí"
$ws.Cells.Item(53,19).Value = 3100
$ws.Cells.Item(53,20).Value = 1

$ws.Cells.Item(54,1).Value = 6
$ws.Cells.Item(54,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(54,3).Value = "Metropolitana"
$ws.Cells.Item(54,4).Value = 44466
$ws.Cells.Item(54,5).Value = 13
$ws.Cells.Item(54,6).Value = "Fruta"
$ws.Cells.Item(54,7).Value = 100107
$ws.Cells.Item(54,8).Value = "Otros"
$ws.Cells.Item(54,9).Value = 100107002
$ws.Cells.Item(54,10).Value = "Chirimoya"
$ws.Cells.Item(54,11).Value = "Cultivar IV Región"
$ws.Cells.Item(54,12).Value = "Extra (doble especial)"
$ws.Cells.Item(54,13).Value = 125
$ws.Cells.Item(54,14).Value = 3400
$ws.Cells.Item(54,15).Value = 3400
$ws.Cells.Item(54,16).Value = 3400
$ws.Cells.Item(54,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(54,18).Value = "Provincia de Limarí"
$ws.Cells.Item(54,19).Value = 3400
$ws.Cells.Item(54,20).Value = 1

$ws.Cells.Item(55,1).Value = 6
$ws.Cells.Item(55,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(55,3).Value = "Metropolitana"
$ws.Cells.Item(55,4).Value = 44466
$ws.Cells.Item(55,5).Value = 13
$ws.Cells.Item(55,6).Value = "Fruta"
$ws.Cells.Item(55,7).Value = 100107
$ws.Cells.Item(55,8).Value = "Otros"
$ws.Cells.Item(55,9).Value = 100107002
$ws.Cells.Item(55,10).Value = "Chirimoya"
$ws.Cells.Item(55,11).Value = "Cultivar IV Región"
$ws.Cells.Item(55,12).Value = "Primera"
$ws.Cells.Item(55,13).Value = 200
$ws.Cells.Item(55,14).Value = 2700
$ws.Cells.Item(55,15).Value = 2700
$ws.Cells.Item(55,16).Value = 2700
$ws.Cells.Item(55,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(55,18).Value = "Provincia de Limarí"
$ws.Cells.Item(55,19).Value = 2700
$ws.Cells.Item(55,20).Value = 1

$ws.Cells.Item(56,1).Value = 6
$ws.Cells.Item(56,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(56,3).Value = "Metropolitana"
$ws.Cells.Item(56,4).Value = 44466
$ws.Cells.Item(56,5).Value = 13
$ws.Cells.Item(56,6).Value = "Fruta"
$ws.Cells.Item(56,7).Value = 100107
$ws.Cells.Item(56,8).Value = "Otros"
$ws.Cells.Item(56,9).Value = 100107002
$ws.Cells.Item(56,10).Value = "Chirimoya"
$ws.Cells.Item(56,11).Value = "Cultivar IV Región"
$ws.Cells.Item(56,12).Value = "Segunda"
$ws.Cells.Item(56,13).Value = 200
$ws.Cells.Item(56,14).Value = 2200
$ws.Cells.Item(56,15).Value = 2200
$ws.Cells.Item(56,16).Value = 2200
$ws.Cells.Item(56,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(56,18).Value = "Provincia de Limarí"
$ws.Cells.Item(56,19).Value = 2200
$ws.Cells.Item(56,20).Value = 1

Write-Output "Inserted rows 53-56; new dimension should be A1:T103"
